$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.212.48'
$ws.Range("E2").Value = '  +2.55%  '
$ws.Range("D3").Value = '3.266.59'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '398.42'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = '108.88'
$ws.Range("E7").Value = '  +4.34%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = '39.41'
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("D11").Value = '0.0956'
$ws.Range("E11").Value = '  +5.40%  '
$ws.Range("E12").Value = '  +1.70%  '
$ws.Range("D13").Value = '3.779.21'
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("D14").Value = '8.28'
$ws.Range("E14").Value = '  +2.29%  '
$ws.Range("D15").Value = '19.01'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").Value = '3.272.42'
$ws.Range("E16").Value = '  +1.42%  '
$ws.Range("E17").Value = '  -1.68%  '
$ws.Range("D18").Value = '11.03'
$ws.Range("E18").Value = '  +3.01%  '
$ws.Range("D19").Value = '57.019.04'
$ws.Range("E19").Value = '  +2.39%  '
$ws.Range("E20").Value = '  -1.03%  '
$ws.Range("E21").Value = '  +5.91%  '
$ws.Range("D22").Value = '12.99'
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("D23").Value = '294.31'
$ws.Range("E23").Value = '  -2.79%  '
$ws.Range("D24").Value = '74.26'
$ws.Range("E24").Value = '  -1.24%  '
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("D26").Value = '28.16'
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '7.93'
$ws.Range("E27").Value = '  -3.37%  '
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '0.112'
$ws.Range("E32").Value = '  +1.00%  '
$ws.Range("D33").Value = '11.22'
$ws.Range("E33").Value = '  -0.35%  '
$ws.Range("D34").Value = '40.19'
$ws.Range("E34").Value = '  +11.27%  '
$ws.Range("D35").Value = '0.0491'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("E36").Value = '  +1.20%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").Value = '3.01'
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("D41").Value = '136.47'
$ws.Range("E41").Value = '  +1.07%  '
$ws.Range("D42").Value = '0.122'
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").Value = '0.286'
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = '3.95'
$ws.Range("E44").Value = '  -1.85%  '
$ws.Range("E45").Value = '  -2.35%  '
$ws.Range("D46").Value = '16.86'
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("D47").Value = '22.44'
$ws.Range("E47").Value = '  +0.65%  '
$ws.Range("E48").Value = '  +5.48%  '
$ws.Range("D49").Value = '2.149.91'
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = '2.39'
$ws.Range("E50").Value = '  -3.70%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '1.98'
$ws.Range("E51").Value = '  -6.36%  '